$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row (2026/02/14, 土, 18, 201) was inserted before the existing
# row 809, shifting the former rows 809-850 down to 810-851 (and the sheet
# dimension grows from D850 to D851).
$ws.Rows("809:809").Insert()

# Column A holds plain text dates (t="inlineStr" in the source), so force
# text interpretation while writing the date-like string, then restore the
# default "Normal" style so the new cell matches its neighbours (which carry
# no explicit style index).
$a = $ws.Range("A809")
$a.NumberFormat = "@"
$a.Value = "2026/02/14"
$a.Style = "Normal"

$ws.Range("B809").Value = "土"
$ws.Range("C809").Value = 18
$ws.Range("D809").Value = 201
